# Update gh-pages output data (attendee counts / prices / one event refresh)
# across the "展览" (sheet1), "演出" (sheet2), "本地生活" (sheet3) and
# "全部类型" (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsAll = $wb.Worksheets.Item("全部类型")

# ----------------------------------------------------------------------
# Sheet "展览" (sheet1)
# ----------------------------------------------------------------------
$wsExpo.Range("F3").Value = 673
$wsExpo.Range("F4").Value = 1485
$wsExpo.Range("F5").Value = 3237
$wsExpo.Range("F7").Value = 657
$wsExpo.Range("F8").Value = 2216
$wsExpo.Range("F9").Value = 477
$wsExpo.Range("F10").Value = 407
$wsExpo.Range("F12").Value = 126
$wsExpo.Range("F13").Value = 307
$wsExpo.Range("F14").Value = 1066
$wsExpo.Range("F17").Value = 76
$wsExpo.Range("F18").Value = 200
$wsExpo.Range("F19").Value = 4433
$wsExpo.Range("F20").Value = 1289
$wsExpo.Range("F21").Value = 3368
$wsExpo.Range("F22").Value = 322
$wsExpo.Range("F23").Value = 70
$wsExpo.Range("F24").Value = 164
$wsExpo.Range("F25").Value = 3303
$wsExpo.Range("G25").Value = 67.5
$wsExpo.Range("F26").Value = 4907
$wsExpo.Range("F29").Value = 540
$wsExpo.Range("F30").Value = 3181
$wsExpo.Range("F31").Value = 345
$wsExpo.Range("F36").Value = 1149
$wsExpo.Range("G36").Value = 99
$wsExpo.Range("F37").Value = 1391
$wsExpo.Range("F38").Value = 113
$wsExpo.Range("F39").Value = 1319
$wsExpo.Range("F40").Value = 841
$wsExpo.Range("F42").Value = 791
$wsExpo.Range("F45").Value = 284

# ----------------------------------------------------------------------
# Sheet "演出" (sheet2)
# ----------------------------------------------------------------------
$wsShow.Range("F7").Value = 995

# ----------------------------------------------------------------------
# Sheet "本地生活" (sheet3)
# ----------------------------------------------------------------------
$wsLocal.Range("F2").Value = 2096

# ----------------------------------------------------------------------
# Sheet "全部类型" (sheet4)
# ----------------------------------------------------------------------
$wsAll.Range("F2").Value = 2096
$wsAll.Range("F3").Value = 673
$wsAll.Range("F4").Value = 1485
$wsAll.Range("F5").Value = 3237
$wsAll.Range("F7").Value = 657
$wsAll.Range("F9").Value = 2216
$wsAll.Range("F10").Value = 477
$wsAll.Range("F11").Value = 407
$wsAll.Range("F13").Value = 995
$wsAll.Range("F14").Value = 126
$wsAll.Range("F15").Value = 307
$wsAll.Range("F16").Value = 1066
$wsAll.Range("F19").Value = 200
$wsAll.Range("F20").Value = 4433
$wsAll.Range("F21").Value = 1289
$wsAll.Range("F23").Value = 3368
$wsAll.Range("F24").Value = 3303
$wsAll.Range("G24").Value = 67.5
$wsAll.Range("F25").Value = 4907
$wsAll.Range("F28").Value = 3181
$wsAll.Range("F29").Value = 345

# Row 34 on "全部类型" is refreshed wholesale to reflect a different event
# (same one now shown, updated, in row 36 of "展览"):
$wsAll.Range("C34").Value = "杭州·夏之誓国乙only-日夜场"
$wsAll.Range("D34").Value = "北干街道萧杭路689号 杭州时尚外滩艺术中心"
$wsAll.Range("E34").Value = "2024.07.27 10:00-07.27 21:00"
$wsAll.Range("F34").Value = 1149
$wsAll.Range("G34").Value = 99
$wsAll.Range("H34").Value = "https://show.bilibili.com/platform/detail.html?id=83589"
$wsAll.Range("I34").Value = "//i2.hdslb.com/bfs/openplatform/202405/99kWb2dy1714964533903.png"

$wsAll.Range("F35").Value = 1391
$wsAll.Range("F36").Value = 113
$wsAll.Range("F37").Value = 1319
$wsAll.Range("F39").Value = 841
$wsAll.Range("F44").Value = 284
